$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Kontakttime" entry: the exam/contact-hour date moved from
# 06.10 to 06.11 ("la inn siste eksamen og oppdaterte kalender").
$ws.Range("C13").Value = " 06.11: Kontakttime, kursansvarlig tilgjengelig i Aud A"

# Update current selection to match the author's final cursor position.
$ws.Range("C14").Select()
